$wb = $excel.ActiveWorkbook

# Add a new worksheet "Contexts" at the end of the workbook (after "Intersections")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "Contexts"

# Populate the new sheet with its header + example data
$data = @(
    @("[Context]", "relations"),
    @("Context",   "Relation"),
    @("Braga",     "t1"),
    @("Braga",     "t2"),
    @("Braga",     "t3")
)
for ($r = 0; $r -lt $data.Length; $r++) {
    $newSheet.Cells.Item($r + 1, 1).Value = $data[$r][0]
    $newSheet.Cells.Item($r + 1, 2).Value = $data[$r][1]
}

# Update the selection on the "Relations" sheet to a whole-column selection (A:B)
$wsRelations = $wb.Worksheets.Item("Relations")
$wsRelations.Activate()
$wsRelations.Range("A1:B1048576").Select()

# Leave the new "Contexts" sheet as the active / selected tab, with B1 selected
$newSheet.Activate()
$newSheet.Range("B1").Select()
